# Replace the "Nur Funktionen und Structs" list item through the
# "Token Klassen" heading with the updated OOXML:
#   - the original list item loses its trailing bookmark
#   - a new empty paragraph, a new "Was nicht geht" heading (berschrift3)
#     and a new "Structs mit Array Feld" list item (which now carries the
#     _GoBack bookmark) are inserted before the page-break paragraph
#   - the stray w:lang="en-US" run/paragraph properties on the page-break
#     paragraph and the "Token Klassen" heading are dropped
$d = $word.ActiveDocument

# Locate the paragraph that ends the span ("Nur Funktionen und Structs")
# and the paragraph that ends the span ("Token Klassen") so the edit does
# not depend on hard-coded character offsets.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($startPara -eq $null -and $t -like "*und Structs*") {
        $startPara = $p
    }
    if ($startPara -ne $null -and $t -like "*Token Klassen*") {
        $endPara = $p
        break
    }
}

$target = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    '<w:p w:rsidR="00892A57" w:rsidRPr="00892A57" w:rsidRDefault="00892A57" w:rsidP="00892A57"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Nur </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Funktionen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> und Structs</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="berschrift3"/></w:pPr><w:r><w:t xml:space="preserve">Was </w:t></w:r><w:r><w:t>nicht geht</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Structs</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> mit Array</w:t></w:r><w:r><w:t xml:space="preserve"> Feld</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' +
    '<w:p w:rsidR="00D204A0" w:rsidRDefault="00D204A0"><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="1F3763" w:themeColor="accent1" w:themeShade="7F"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:br w:type="page"/></w:r></w:p>' +
    '<w:p w:rsidR="00D204A0" w:rsidRDefault="00D204A0" w:rsidP="00D204A0"><w:pPr><w:pStyle w:val="berschrift5"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Token Klassen</w:t></w:r></w:p>' +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xmlFrag)
